$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 13 rows (146-157) no longer carry the team-icon Image column value
$ws.Range("E146:E157").ClearContents()

# Week 14 results (rows 158-169)
$week14 = @(
    @("Kauaireek Hill", 91.16, 76.56, "Team Icons/kauaireek-modified.png"),
    @("Chasing dank Herb", 100.08, 107.68, "Team Icons/chasing-modified.png"),
    @("Ju Ju Smith Poopster", 85.92, 101.18, "Team Icons/juju-modified.png"),
    @("Bye Breece See You in ValHalla", 101.18, 85.92, "Team Icons/breece-modified.png"),
    @("Cooking with Gas", 76.56, 91.16, "Team Icons/cooking-modified.png"),
    @("Dulcich de Leche", 109.76, 117.02, "Team Icons/dulcich-modified.png"),
    @("Dillon Panthers", 110.88, 105.6, "Team Icons/dillon-modified.png"),
    @("Daemon and the Rightful Heirs", 89.62, 115.26, "Team Icons/daemon-modified.png"),
    @("Krombopulos Michael Evans", 117.02, 109.76, "Team Icons/krombopulos-modified.png"),
    @("Freier Freier Pants on Fire", 107.68, 100.08, "Team Icons/freier-modified.png"),
    @("Christian Kirk Cousins", 105.6, 110.88, "Team Icons/sir-modified.png"),
    @("Fantasy Football Champion 2022", 115.26, 89.62, "Team Icons/fantasy-modified.png")
)

$startRow = 158
for ($i = 0; $i -lt $week14.Length; $i++) {
    $row = $startRow + $i
    # Inserting a row just below the existing Week-13 block carries the
    # C:D number style (s="1") down from the row above, matching how the
    # sheet already looks for every prior week's block.
    $ws.Rows($row).Insert()
}
for ($i = 0; $i -lt $week14.Length; $i++) {
    $row = $startRow + $i
    $data = $week14[$i]
    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = "Week 14"
    $ws.Cells.Item($row, 3).Value = $data[1]
    $ws.Cells.Item($row, 4).Value = $data[2]
    $ws.Cells.Item($row, 5).Value = $data[3]
}

# Match the author's final selection/scroll position (one row past the
# new last data row, in column D) as recorded in the saved view state.
[void]($excel.ActiveWindow.ScrollRow = 143)
[void]$ws.Range("D171").Select()
